# Refresh the cryptos price list (GitHub Actions scheduled update).
# Column D ("Price") values are numeric-looking strings (e.g. "1.00", "28.01")
# that must stay as literal text -- Excel auto-converts a plain assignment of
# such a string to a real number (dropping trailing zeros / introducing float
# rounding, e.g. "40.62" -> 40.619999999999997). Forcing NumberFormat="@" (Text)
# before the write keeps it text, then Style="Normal" resets formatting so no
# stray style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.347.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.646.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.93%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.645.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("E10").Value = "  +8.05%  "

$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("E13").Value = "  +2.55%  "

$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000193"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.27%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.126.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.310.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.644.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.64%  "

$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "

$ws.Range("E23").Value = "  +2.42%  "

$ws.Range("E24").Value = "  +0.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "75.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("E28").Value = "  +2.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.781.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "561.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.98%  "

$ws.Range("E32").Value = "  +2.55%  "

$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("E35").Value = "  +2.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("E37").Value = "  +5.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.376"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.81%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0339"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "

$ws.Range("E44").Value = "  +0.94%  "

$ws.Range("E45").Value = "  +2.13%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.67%  "

$ws.Range("E49").Value = "  +2.29%  "

$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.47%  "
